$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 2D training schedule data (rows 2-6, columns A-J)
$data = @(
    @(1, 5, 3, 1, 5, -4, 2, 54, 5, "train_dim2_1"),
    @(2, 6, 2, 1, 3, -5, 1, 65, 5, "train_dim2_1"),
    @(3, 6, 4, 5, 9, -1, 5, 21, 5, "train_dim2_1"),
    @(4, 5, 1, 2, 4, -3, 3, 43, 5, "train_dim2_1"),
    @(5, 8, 4, 6, 8, -2, 4, 32, 5, "train_dim2_1")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}

# Update selected cell to match new active selection
$ws.Range("I1").Select()

$wb.Save()
